# The document has two logos repeated in its headers/footers:
#   - the BTEC logo (a .jpg) living in both page headers
#   - the Pearson Edexcel logo (a .png) living in both page footers
#
# The commit swaps the "name" each picture's drawing object reports:
#   header pictures : image2.jpg -> image1.jpg
#   footer pictures : image1.png -> image2.png
#
# In the Word object model this is the InlineShape.Name property
# (backed by <wp:docPr name="...">). Drive it through Selection rather
# than the InlineShape reference directly - selecting first keeps the
# COM handle fresh, which matters for shapes that live in footer
# stories.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Set-LogoName($range, [string]$newName) {
    $range.InlineShapes(1).Select() | Out-Null
    $word.Selection.InlineShapes(1).Name = $newName
}

# Headers: BTEC logo, image2.jpg -> image1.jpg
foreach ($hdr in $sec.Headers) {
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        Set-LogoName $hdr.Range "image1.jpg"
    }
}

# Footers: Pearson Edexcel logo, image1.png -> image2.png
foreach ($ftr in $sec.Footers) {
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        Set-LogoName $ftr.Range "image2.png"
    }
}
